$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells keep their original string representation
# (these columns hold formatted text like "27.032.96" or "  +2.24%  ", not numbers)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.032.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.657.87'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.72%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.49'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.42%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.55%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.14'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0884'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.892.42'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.656.04'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.54%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.03%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.523'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.65%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.024.63'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.25%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0738'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.77'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.66%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.43'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.66%  '
$ws.Range("B23").Value = 'Toncoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.23'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.10%  '
$ws.Range("B24").Value = 'Avalanche'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.29'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.08'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.13'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.92%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.86'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.44%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.35%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.557.54'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.30'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.01%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.62'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +7.74%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.22%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.579'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.17%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +8.69%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.55%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.42'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +8.27%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.25'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.80%  '
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.973'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.802.43'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.87%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.18'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.33%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.77%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1000'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.15%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.99%  '
